$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Remove the two knockouts that are now handled in iFerment assignments ---
$ws.Range("C36").Value = 0
$ws.Range("C167").Value = 0

# --- 2. Append the six new reactions that used to be knocked out ---
$newReactions = @("ICDHyr", "ME2", "G3PD2", "HACD1a", "MTHFR3_1", "LDH_L")

$firstRow = 240

for ($i = 0; $i -lt $newReactions.Count; $i++) {
    $row = $firstRow + $i
    $nameCell = $ws.Range("A$row")
    $nameCell.Value = $newReactions[$i]

    if ($i -eq 0) {
        # Build the pasted-in-from-elsewhere look (small orange monospace font,
        # vertically centered) on the very first new cell ...
        $nameCell.Font.Name = "Consolas"
        $nameCell.Font.Size = 8
        $nameCell.Font.Color = 7901646
        $nameCell.VerticalAlignment = -4108
    } else {
        # ... then reuse that exact formatting for the rest via copy/paste so we
        # don't keep minting new styles.
        $ws.Range("A$firstRow").Copy() | Out-Null
        $ws.Range("A$row").PasteSpecial(-4122) | Out-Null
        $ws.Range("A$row").Value = $newReactions[$i]
    }

    # Match the row height used throughout the rest of the sheet.
    $ws.Rows.Item($row).RowHeight = 18

    # Columns B-F: copy the existing formatting used by every other data row
    # (reuses the existing styles instead of creating new ones) and zero them out.
    $ws.Range("B239:F239").Copy() | Out-Null
    $ws.Range("B$row`:F$row").PasteSpecial(-4122) | Out-Null
    $ws.Range("B$row`:F$row").Value = 0
}

$excel.CutCopyMode = 0

# --- 3. Leave the selection where the author's last edit was ---
$ws.Range("C36").Select() | Out-Null

Write-Output "done"
